$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / tab
$ws.Name = "Through 2022-09-25"

# Update header label for the September row
$ws.Range("A10").Value = "September (through 09-25)"

# Update September figures (row 10)
$ws.Range("C10").Value = 36
$ws.Range("D10").Value = 64
$ws.Range("E10").Value = 49
$ws.Range("F10").Value = 61
$ws.Range("G10").Value = 98
$ws.Range("H10").Value = 152
$ws.Range("I10").Value = 120

# Update Total figures (row 11)
$ws.Range("C11").Value = 417
$ws.Range("D11").Value = 615
$ws.Range("E11").Value = 539
$ws.Range("F11").Value = 410
$ws.Range("G11").Value = 882
$ws.Range("H11").Value = 1222
$ws.Range("I11").Value = 1255
